# The edit swaps the data of rows 16 and 17 on the active sheet (two
# species-observation records traded places), while every other row is
# left untouched.
#
# We swap cell-by-cell (only the columns whose content actually differs
# between the two rows) using Range.Copy through a scratch cell on row
# 100 as a temporary holding spot, then clear the scratch cell. Copy
# (unlike assigning .Value/.Value2, which re-parses the payload the way
# typing into the Excel UI would and can turn numeric-looking text such
# as "25" into a real number, or "2023-08-29" into a date) preserves the
# exact stored type of each source cell, so text stays text and numbers
# stay numbers with full precision.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","D","E","F","G","H","I","J","Q","R","Z","AB")

foreach ($col in $cols) {
    $src16 = $ws.Range($col + "16")
    $src17 = $ws.Range($col + "17")
    $scratch = $ws.Range($col + "100")

    $src16.Copy($scratch)
    $src17.Copy($src16)
    $scratch.Copy($src17)
    $scratch.Clear()
}

# Column L ("Kön") holds an empty placeholder cell that also trades
# places: row 17 had it, row 16 did not. After the swap row 17 no
# longer carries that placeholder.
$ws.Range("L17").Value = ""
